# Created test case for Read More Project
# - Simplifies the "Outcome" column (G) for most already-passing test cases
#   to a standard "Passed test case" label.
# - Updates the CreateProject test case's outcome wording.
# - Adds a brand new test case row (row 16) for "Read More" on the
#   portfolio project listing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Table 1 (Portfolio Website test scenarios): simplify Outcome (col G) ---
$passed = "Passed test case"
$ws.Cells.Item(3, 7).Value  = $passed
$ws.Cells.Item(4, 7).Value  = $passed
$ws.Cells.Item(5, 7).Value  = $passed
$ws.Cells.Item(6, 7).Value  = $passed
$ws.Cells.Item(7, 7).Value  = $passed
$ws.Cells.Item(8, 7).Value  = $passed
$ws.Cells.Item(9, 7).Value  = $passed
$ws.Cells.Item(10, 7).Value = $passed
$ws.Cells.Item(11, 7).Value = $passed
$ws.Cells.Item(12, 7).Value = $passed

# Row 15 (test_<CreateProject>) outcome re-worded
$ws.Cells.Item(15, 7).Value = "Not able to test this function, but based on the given source code, project object can be created through command prompt"

# --- New row 16: test_<ReadMorePortfolio> ---
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "test_<ReadMorePortfolio>"
$ws.Cells.Item(16, 3).Value = "This is to test whether users can click on 'Read More' to display project details"
$ws.Cells.Item(16, 4).Value = "NIL"
$ws.Cells.Item(16, 5).Value = "Project details displayed"
$ws.Cells.Item(16, 7).Value = "Failed test case"

# --- Table 2 (Portfolio Admin Site test scenarios): simplify Outcome (col G) ---
$ws.Cells.Item(23, 7).Value = $passed
$ws.Cells.Item(24, 7).Value = $passed
$ws.Cells.Item(25, 7).Value = $passed
$ws.Cells.Item(26, 7).Value = $passed
$ws.Cells.Item(27, 7).Value = $passed
$ws.Cells.Item(28, 7).Value = $passed
$ws.Cells.Item(29, 7).Value = $passed

# --- Restore the view to the top of the sheet and select the new row's
#     outcome cell, mirroring where the author ended up after typing in the
#     new test case. ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 4
$null = $ws.Range("H16").Select()
